$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.008.94'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '3.396.18'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.72%  '
$ws.Range('D6').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.05%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.397.21'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('D9').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.540'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.96%  '
$ws.Range('D10').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('D11').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.121'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.31%  '
$ws.Range('D12').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.432'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.70%  '
$ws.Range('D13').Value = '3.986.82'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('D14').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.133'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.75%  '
$ws.Range('E15').Value = '  +3.75%  '
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').Value = '63.142.96'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '3.404.45'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D20').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('D21').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.43%  '
$ws.Range('D22').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.35%  '
$ws.Range('D23').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.995'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('D24').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.50'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.41%  '
$ws.Range('E25').Value = '  -2.22%  '
$ws.Range('E26').Value = '  +20.58%  '
$ws.Range('D27').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.44%  '
$ws.Range('E28').Value = '  -2.39%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.97'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.91%  '
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('D33').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.18%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.63'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.85%  '
$ws.Range('E38').Value = '  -2.05%  '
$ws.Range('D39').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0758'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.28%  '
$ws.Range('D40').Value = '2.902.95'
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('E41').Value = '  -3.74%  '
$ws.Range('D42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.66'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.74%  '
$ws.Range('E43').Value = '  +1.89%  '
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.23%  '
$ws.Range('E48').Value = '  +1.67%  '
$ws.Range('E49').Value = '  +18.89%  '
$ws.Range('D50').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.37'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.90%  '
$ws.Range('E51').Value = '  +2.99%  '
